$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix: document number column (G) parsed as text when it actually holds
# a pure numeric value - convert to real numbers.
$ws.Range("G4").Value = 124254
$ws.Range("G5").Value = 15005

# Document Identity: the "A" column now carries the user identity
# (login) that touched the row instead of the old running row number.
$ws.Range("A5").Value = "vasia02"
$ws.Range("A6").Value = "petia34"
$ws.Range("A7").Clear()
$ws.Range("A8").Value = "oleg"

$ws.Range("A8").Select()
